$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "TEST_IMPORT_SURVEY_RESP_1_test"
